# Automatische test-sync: 2025-06-19 21:43:50
# Add new row 26 to the "Logs" sheet describing an incoming mail, then
# update the dependent ranges/counters (conditional formatting + Dashboard tally).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A26").Value = "Uitnodiging voor netwerkevent"
$ws.Range("B26").Value = "mailmind.test@zohomail.eu"
$ws.Range("C26").Value = "Graag nodig ik u uit voor ons zakelijke netwerkevent volgende maand."
$ws.Range("D26").Value = "Samenwerking / Partnerverzoek"
$ws.Range("F26").Value = "2025-06-19 21:43:13"
$ws.Range("G26").Value = "Nee"

# Extend the conditional formatting ranges to cover the new row.
$catFcs = $ws.Range("D2:D25").FormatConditions
for ($i = 1; $i -le $catFcs.Count; $i++) {
    $catFcs.Item($i).ModifyAppliesToRange($ws.Range("D2:D26"))
}

$answeredFcs = $ws.Range("G2:G25").FormatConditions
for ($i = 1; $i -le $answeredFcs.Count; $i++) {
    $answeredFcs.Item($i).ModifyAppliesToRange($ws.Range("G2:G26"))
}

# Bump the "Samenwerking / Partnerverzoek" tally on the Dashboard sheet.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 7
